$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 547 (shifts existing rows 547-624 down to 548-625)
$ws.Rows.Item(547).Insert()

$ws.Cells.Item(547, 1).Value = 5
$ws.Cells.Item(547, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(547, 3).Value = "Maule"
$ws.Cells.Item(547, 4).Value = 45127
$ws.Cells.Item(547, 5).Value = 7
$ws.Cells.Item(547, 6).Value = 100114014
$ws.Cells.Item(547, 7).Value = "Betarraga"
$ws.Cells.Item(547, 8).Value = "Sin especificar"
$ws.Cells.Item(547, 9).Value = "Primera"
$ws.Cells.Item(547, 10).Value = 5000
$ws.Cells.Item(547, 11).Value = 550
$ws.Cells.Item(547, 12).Value = 550
$ws.Cells.Item(547, 13).Value = 550
$ws.Cells.Item(547, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(547, 15).Value = "Región del Maule"
$ws.Cells.Item(547, 16).Value = 110
$ws.Cells.Item(547, 17).Value = 5
$ws.Cells.Item(547, 18).Value = "Hortaliza"
